# Updates the cryptos list (price / 1h volume columns, and a few
# re-ranked rows) to match the refreshed GitHub Actions data pull.
# Numeric-looking "Price" strings are forced to remain text (NumberFormat
# "@" then Style back to Normal) so Excel doesn't silently coerce them
# into real numbers, matching the original text cells in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '73.305.50'
$ws.Range('E2').Value = '  +1.88%  '

$ws.Range('D3').Value = '4.052.15'
$ws.Range('E3').Value = '  +1.09%  '

$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.12'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +9.09%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.76'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.29%  '

$ws.Range('D7').Value = '4.046.11'
$ws.Range('E7').Value = '  +1.04%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.692'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.50%  '

$ws.Range('E9').Value = '  -0.01%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.762'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.16%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.170'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.28%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.84'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +13.02%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000326'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.07%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '11.17'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.35%  '

$ws.Range('D15').Value = '4.697.25'
$ws.Range('E15').Value = '  +0.98%  '

$ws.Range('D16').Value = '4.053.97'
$ws.Range('E16').Value = '  +1.06%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.33'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.59%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '20.85'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.66%  '

$ws.Range('E19').Value = '  +3.19%  '

$ws.Range('B20').Value = 'WrappedBTC'
$ws.Range('C20').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D20').Value = '73.217.82'
$ws.Range('E20').Value = '  +1.89%  '

$ws.Range('B21').Value = 'TRON'
$ws.Range('C21').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.132'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.21%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '444.74'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.42%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '98.13'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.36%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.54'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +8.15%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.56'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.89%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '14.60'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.72%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.25'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +19.09%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.41'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.28%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '11.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.65%  '

$ws.Range('E30').Value = '  +2.12%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '37.07'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.05%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.86'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +10.66%  '

$ws.Range('E33').Value = '  +3.96%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '13.70'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.63%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '691.81'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.38%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '48.51'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +10.20%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '68.28'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.87%  '

$ws.Range('E38').Value = '  +9.49%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.447'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.89%  '

$ws.Range('E40').Value = '  -2.04%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '11.44'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +18.82%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.38'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.02%  '

$ws.Range('E43').Value = '  +0.06%  '

$ws.Range('E44').Value = '  +2.55%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0495'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.69%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.999'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.06%  '

$ws.Range('E47').Value = '  +0.88%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.76'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.07%  '

$ws.Range('B49').Value = 'LidoDAOToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.50'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.52%  '

$ws.Range('B50').Value = 'ARBITRUM'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.22'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +11.27%  '

$ws.Range('B51').Value = 'ApeXProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.32'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.14%  '
